$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header ("feat" / "shap") is unchanged.
# Rows 2-70: feature-index (A), feature-name (B) and SHAP value (C)
# re-sorted ascending by the (recomputed) SHAP value.
$rowData = @(
    @{Row=2; A=2; B="coulomb"; C="0"},
    @{Row=3; A=30; B="occr4"; C="0"},
    @{Row=4; A=29; B="occr3"; C="0"},
    @{Row=5; A=28; B="occr2"; C="0"},
    @{Row=6; A=27; B="occr1"; C="0"},
    @{Row=7; A=46; B="SCFOCCr4"; C="0"},
    @{Row=8; A=45; B="SCFOCCr3"; C="0"},
    @{Row=9; A=44; B="SCFOCCr2"; C="0"},
    @{Row=10; A=43; B="SCFOCCr1"; C="0"},
    @{Row=11; A=65; B="SCFFp"; C="0"},
    @{Row=12; A=66; B="SCFFq"; C="0"},
    @{Row=13; A=1; B="pair_energy"; C="7.716963870198247e-11"},
    @{Row=14; A=50; B="SCFOCCs4"; C="2.325172841517294e-07"},
    @{Row=15; A=48; B="SCFOCCs2"; C="2.342662069568913e-07"},
    @{Row=16; A=49; B="SCFOCCs3"; C="2.426294550217148e-07"},
    @{Row=17; A=47; B="SCFOCCs1"; C="2.777837984205288e-07"},
    @{Row=18; A=60; B="hqq"; C="2.884392512463198e-07"},
    @{Row=19; A=68; B="SCFOCCq"; C="2.887254040271499e-07"},
    @{Row=20; A=67; B="SCFOCCp"; C="2.991790885353124e-07"},
    @{Row=21; A=31; B="occs1"; C="8.583691881192582e-07"},
    @{Row=22; A=33; B="occs3"; C="1.053489669616333e-06"},
    @{Row=23; A=32; B="occs2"; C="1.083748446022736e-06"},
    @{Row=24; A=0; B="From_Same_Orbital"; C="1.885888150651915e-06"},
    @{Row=25; A=34; B="occs4"; C="3.677954911480495e-06"},
    @{Row=26; A=38; B="SCFFr4"; C="4.166896442014128e-06"},
    @{Row=27; A=14; B="eijab_4"; C="6.082808949942495e-06"},
    @{Row=28; A=35; B="SCFFr1"; C="6.267799861822713e-06"},
    @{Row=29; A=63; B="occp"; C="7.644434890804024e-06"},
    @{Row=30; A=37; B="SCFFr3"; C="9.786079668646735e-06"},
    @{Row=31; A=64; B="occq"; C="1.160236817893111e-05"},
    @{Row=32; A=15; B="screenvirt_1"; C="1.246321771416641e-05"},
    @{Row=33; A=9; B="screen2_3"; C="1.285734420393612e-05"},
    @{Row=34; A=17; B="screenvirt_3"; C="1.533054277918431e-05"},
    @{Row=35; A=40; B="SCFFs2"; C="1.534470392386801e-05"},
    @{Row=36; A=56; B="hss2"; C="1.67000580812405e-05"},
    @{Row=37; A=55; B="hss1"; C="1.689288264658724e-05"},
    @{Row=38; A=13; B="eijab_3"; C="1.809746288411477e-05"},
    @{Row=39; A=58; B="hss4"; C="1.812935385762642e-05"},
    @{Row=40; A=25; B="Fs3"; C="1.824138155036271e-05"},
    @{Row=41; A=10; B="screen2_4"; C="1.92359541372374e-05"},
    @{Row=42; A=54; B="hrr4"; C="2.118160951747545e-05"},
    @{Row=43; A=6; B="screen1_4"; C="2.167753923753714e-05"},
    @{Row=44; A=62; B="Fq"; C="2.393468630956089e-05"},
    @{Row=45; A=36; B="SCFFr2"; C="2.479385293874846e-05"},
    @{Row=46; A=21; B="Fr3"; C="2.534631068421036e-05"},
    @{Row=47; A=23; B="Fs1"; C="2.7245629097016e-05"},
    @{Row=48; A=53; B="hrr3"; C="2.834766055828733e-05"},
    @{Row=49; A=12; B="eijab_2"; C="3.04776222587208e-05"},
    @{Row=50; A=11; B="eijab_1"; C="3.135206226056221e-05"},
    @{Row=51; A=42; B="SCFFs4"; C="3.178245358218824e-05"},
    @{Row=52; A=18; B="screenvirt_4"; C="3.207539147244977e-05"},
    @{Row=53; A=51; B="hrr1"; C="3.288769948605978e-05"},
    @{Row=54; A=16; B="screenvirt_2"; C="3.494307314350983e-05"},
    @{Row=55; A=41; B="SCFFs3"; C="3.514199163176058e-05"},
    @{Row=56; A=52; B="hrr2"; C="3.650377376256499e-05"},
    @{Row=57; A=22; B="Fr4"; C="3.866612528759924e-05"},
    @{Row=58; A=5; B="screen1_3"; C="4.151725316181081e-05"},
    @{Row=59; A=26; B="Fs4"; C="4.245195593576835e-05"},
    @{Row=60; A=39; B="SCFFs1"; C="4.522295191001501e-05"},
    @{Row=61; A=3; B="screen1_1"; C="4.823712106852731e-05"},
    @{Row=62; A=59; B="hpp"; C="4.980231149011689e-05"},
    @{Row=63; A=8; B="screen2_2"; C="5.322045189521259e-05"},
    @{Row=64; A=20; B="Fr2"; C="5.352977670267405e-05"},
    @{Row=65; A=4; B="screen1_2"; C="6.351828417088331e-05"},
    @{Row=66; A=19; B="Fr1"; C="6.541616762866939e-05"},
    @{Row=67; A=24; B="Fs2"; C="6.593858394030188e-05"},
    @{Row=68; A=61; B="Fp"; C="6.915529427541604e-05"},
    @{Row=69; A=57; B="hss3"; C="7.398755716951347e-05"},
    @{Row=70; A=7; B="screen2_1"; C="0.000180491434114067"}
)

foreach ($r in $rowData) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = [double]$r.C
}
